# 🔄 Actualización automática del tracker
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update resultado/profit for several rows
$ws.Range("G38").Value = "Acierto"
$ws.Range("H38").Value = 1.75

$ws.Range("G39").Value = "Fallo"
$ws.Range("H39").Value = -1

$ws.Range("G43").Value = "Fallo"
$ws.Range("H43").Value = -1

$ws.Range("G46").Value = "Fallo"
$ws.Range("H46").Value = -1

$ws.Range("G52").Value = "Fallo"
$ws.Range("H52").Value = -1

$ws.Range("G53").Value = "Fallo"
$ws.Range("H53").Value = -1

$ws.Range("G54").Value = "Fallo"
$ws.Range("H54").Value = -1

$ws.Range("G56").Value = "Fallo"
$ws.Range("H56").Value = -1

# Fix event_id type: was stored as text, should be numeric
$ws.Range("A63").Value = 14579389
$ws.Range("A64").Value = 14579388
